# Capital calls: Call amount is now sum of remittance call amounts.
# This adds an "Exchange Rates" sheet and re-points the CapitalCall
# Call Date / Due Date columns at the (now per-remittance) dates.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update CapitalCall (sheet1) dates ---------------------------------
# Row 2: both Call Date and Due Date -> 10 Jan 2022 (44571)
$ws1.Range("E2").Value = 44571
$ws1.Range("F2").Value = 44571

# Row 3: both Call Date and Due Date -> 10 Jun 2022 (44722)
$ws1.Range("E3").Value = 44722
$ws1.Range("F3").Value = 44722

# Row 4: both Call Date and Due Date -> 10 Dec 2022 (44905, unchanged value,
# but Due Date now matches Call Date's format/value)
$ws1.Range("E4").Value = 44905
$ws1.Range("F4").Value = 44905

# Both date columns now share the same dd/mm/yy format.
$ws1.Range("E2:F4").NumberFormat = "dd/mm/yy"

# --- Add the new "Exchange Rates" sheet, right after CapitalCall -------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Exchange Rates"

$ws2.Range("A1").Value = "From Currency"
$ws2.Range("B1").Value = "To Currency"
$ws2.Range("C1").Value = "Exchange Rate "
$ws2.Range("D1").Value = "As Of"

$ws2.Range("A2").Value = "USD"
$ws2.Range("B2").Value = "INR"
$ws2.Range("C2").Value = 80
$ws2.Range("D2").Value = 44571

$ws2.Range("A3").Value = "USD"
$ws2.Range("B3").Value = "INR"
$ws2.Range("C3").Value = 81
$ws2.Range("D3").Value = 44722

$ws2.Range("A4").Value = "USD"
$ws2.Range("B4").Value = "INR"
$ws2.Range("C4").Value = 82
$ws2.Range("D4").Value = 44905

$ws2.Range("A5").Value = "SGD"
$ws2.Range("B5").Value = "INR"
$ws2.Range("C5").Value = 80
$ws2.Range("D5").Value = 44571

$ws2.Range("A6").Value = "SGD"
$ws2.Range("B6").Value = "INR"
$ws2.Range("C6").Value = 81
$ws2.Range("D6").Value = 44722

$ws2.Range("A7").Value = "SGD"
$ws2.Range("B7").Value = "INR"
$ws2.Range("C7").Value = 82
$ws2.Range("D7").Value = 44905

$ws2.Range("D2:D7").NumberFormat = "dd/mm/yy"

# --- Restore selections: Exchange Rates keeps D2:D4 selected, but the
# CapitalCall tab stays the active/selected one with F2:F4 selected. ----
[void]$ws2.Range("D2:D4").Select()
[void]$ws1.Range("F2:F4").Select()
